$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the specific data cells that were removed in this revision
# (various data bits for usmr dapr3)
$cellsToClear = @(
    "G5",
    "B7", "T7",
    "C8", "N8",
    "L10",
    "F11",
    "P12",
    "I13",
    "H14",
    "D16", "F16",
    "O18",
    "U19",
    "T21",
    "D22", "M22",
    "G23", "K23",
    "H24",
    "U27",
    "C28",
    "F29"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# Update the selected / active cell on the sheet
$ws.Range("C8").Select()
